{"js": "// The document contains several \"<id>...</id>\" markers that were each split\n// across three separate runs (an opening-tag run in Courier New / 7f6000,\n// a plain-black run holding just the id value, and a closing-tag run back\n// in Courier New / 7f6000). Re-downloading these ids collapsed each triple\n// back into a single Courier-New run holding the whole \"<id>...</id>\" text.\n// One of them (p169_2) was also renamed to p169v_2 as part of the refresh.\nconst idUpdates = [\n  { oldId: \"p169_2\", newId: \"p169v_2\" },\n  { oldId: \"p170r_1\", newId: \"p170r_1\" },\n  { oldId: \"p170r_2\", newId: \"p170r_2\" },\n  { oldId: \"p170r_3\", newId: \"p170r_3\" },\n  { oldId: \"p170r_4\", newId: \"p170r_4\" },\n  { oldId: \"p170r_5\", newId: \"p170r_5\" },\n  { oldId: \"p170r_6\", newId: \"p170r_6\" },\n];\n\nfor (const { oldId, newId } of idUpdates) {\n  const searchResults = context.document.body.search(`<id>${oldId}</id>`, {\n    matchCase: true,\n  });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  if (searchResults.items.length === 0) {\n    continue;\n  }\n\n  // Replacing the whole matched range with plain text collapses the three\n  // original runs into one, inheriting the formatting of the run at the\n  // start of the range (the Courier-New \"<id>\" run).\n  searchResults.items[0].insertText(`<id>${newId}</id>`, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document contains several \"<id>...</id>\" markers that were each split\n# across three separate runs (an opening-tag run in Courier New / 7f6000,\n# a plain-black run holding just the id value, and a closing-tag run back\n# in Courier New / 7f6000). Re-downloading these ids collapsed each triple\n# back into a single Courier-New run holding the whole \"<id>...</id>\" text.\n# One of them (p169_2) was also renamed to p169v_2 as part of the refresh.\n$d = $word.ActiveDocument\n\n$updates = @(\n    @{ OldId = \"p169_2\";  NewId = \"p169v_2\" },\n    @{ OldId = \"p170r_1\"; NewId = \"p170r_1\" },\n    @{ OldId = \"p170r_2\"; NewId = \"p170r_2\" },\n    @{ OldId = \"p170r_3\"; NewId = \"p170r_3\" },\n    @{ OldId = \"p170r_4\"; NewId = \"p170r_4\" },\n    @{ OldId = \"p170r_5\"; NewId = \"p170r_5\" },\n    @{ OldId = \"p170r_6\"; NewId = \"p170r_6\" }\n)\n\nforeach ($u in $updates) {\n    $findText = \"<id>$($u.OldId)</id>\"\n    $replaceText = \"<id>$($u.NewId)</id>\"\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    # Use the full Execute overload (with ReplaceWith/Replace) so the runs\n    # spanning the match get merged into one even when the replacement text\n    # is identical to the text that was found (plain \"$range.Text = ...\" is\n    # a no-op in that case).\n    $found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
